$d = $word.ActiveDocument

# Chr(11) is a manual line break (Shift+Enter) character; Word serializes it
# as a <w:br/> element and splits the surrounding text into separate <w:t> runs.
$vt = [char]11

# --- Portuguese paragraph under "Programa" -------------------------------
$ptRange = $d.Content
$ptRange.Find.ClearFormatting()
$found = $ptRange.Find.Execute(
    "A carreira de Engenharia Física. Cientistas x engenheiros*todo.",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $ptRange.Text = (
        "A carreira de Engenharia Física. Cientistas x engenheiros: o papel interdisciplinar da Engenharia Física. Campos de atuação. " + $vt +
        "A Física como ciência conceitual: Como aprender Física. Realização de demonstrações e experimentos científicos significativos de Física." + $vt +
        "Conceitos básicos de Engenharia. Habilidades e competências de um engenheiro." + $vt +
        "Desenvolvimento de um projeto temático de Engenharia Física." + $vt +
        "Competição entre projetos de diferentes grupos." + $vt +
        "Avaliação das competições e da disciplina como um todo."
    )
}

# --- English (italic) paragraph under "Programa" -------------------------
$enRange = $d.Content
$enRange.Find.ClearFormatting()
$found = $enRange.Find.Execute(
    "The career of Engineering Physics. Scientists x engineers*whole.",
    $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $enRange.Text = (
        "The career of Engineering Physics. Scientists x engineers: the interdisciplinary role of Engineering Physics. Fields of action." + $vt +
        "Physics as a conceptual science: How to learn Physics. Realization of demonstrations and significant scientific experiments in Physics." + $vt +
        "Basic engineering concepts. Skills and competences of an engineer." + $vt +
        "Development of a thematic project of Physical Engineering." + $vt +
        "Competition between projects from different groups." + $vt +
        "Evaluation of competitions and the discipline as a whole."
    )
}
